$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 6) down onto the
# two new rows, matching the borders/style used by the other data rows.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F8").PasteSpecial(-4122)

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "PT010"
$ws.Range("C7").Value = "Lê Thị Ngọc Ánh"
$ws.Range("D7").Value = "Ngôn ngữ lập trình C#"
$ws.Range("E7").Value = 23
$ws.Range("F7").Value = 46000

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "PT011"
$ws.Range("C8").Value = "Trần Lê Tuyết Mai"
$ws.Range("D8").Value = "Đại số tuyến tính"
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 66000
